$d = $word.ActiveDocument

# Locate the start of the paragraph beginning with "Aby problem nadawal sie..."
$startRng = $d.Content.Duplicate
$startRng.Find.ClearFormatting()
$startRng.Find.Text = "Aby problem nadawa"
$found1 = $startRng.Find.Execute()
$startRng.Expand(4)

# Locate the end of the paragraph ending with "...jest to poniekad rozwiazywanie..."
$endRng = $d.Content.Duplicate
$endRng.Find.ClearFormatting()
$endRng.Find.Text = "jest to poniek"
$found2 = $endRng.Find.Execute()
$endRng.Expand(4)

$target = $d.Range($startRng.Start, $endRng.End)

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Stylakapitu"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Calibri" w:cs="" w:cstheme="minorBidi" w:eastAsiaTheme="minorHAnsi"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="pl-PL" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="" w:cstheme="minorBidi" w:eastAsiaTheme="minorHAnsi"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="pl-PL" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>Nawet pobieżna analiza problemów, do rozwiązywania których używane są algorytmy genetyczne, pozwala na wyróżnienie kilku cech wspólnych, które sprzyjają a może nawet są konieczne [by dało się zastosować alg.gen.]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Stylakapitu"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="" w:cstheme="minorBidi" w:eastAsiaTheme="minorHAnsi"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="pl-PL" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>[o funkcji celu]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Stylakapitu"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="" w:cstheme="minorBidi" w:eastAsiaTheme="minorHAnsi"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="pl-PL" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>Wiele z nich opisane jest funkcją matematyczną, która może posłużyć za funkcję celu [dla alg.gen.]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Stylakapitu"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="" w:cstheme="minorBidi" w:eastAsiaTheme="minorHAnsi"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="pl-PL" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>[przykład] zadania optymalizacyjne – poszukiwanie ekstremum wspomnianej funkcji</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Stylakapitu"/><w:rPr></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Calibri" w:cs="" w:cstheme="minorBidi" w:eastAsiaTheme="minorHAnsi"/><w:color w:val="auto"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="pl-PL" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr><w:t>[cechy funkcji celu] Pozwala ona na jednoznaczne przyporządkowanie wartości osobnika danej kombinacji jego genów. [Jednoznaczność wynika ze „statyczności” problemu]</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xmlFrag)
